# Update transition-probability matrix values on the active sheet
# (recomputed probabilities after adding more simulated games).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1602564102564103
$ws.Range("C2").Value = 0.6121794871794872
$ws.Range("J2").Value = 0.01602564102564102
$ws.Range("P2").Value = 0.1025641025641026
$ws.Range("S2").Value = 0.108974358974359
$ws.Range("B3").Value = 0.005128205128205128
$ws.Range("C3").Value = 0.01025641025641026
$ws.Range("J3").Value = 0.01538461538461539
$ws.Range("P3").Value = 0.7487179487179487
$ws.Range("S3").Value = 0.2205128205128205
$ws.Range("J4").Value = 0.04081632653061224
$ws.Range("P4").Value = 0.7142857142857143
$ws.Range("S4").Value = 0.2448979591836735
$ws.Range("B6").Value = 0.04624277456647399
$ws.Range("D6").Value = 0.005780346820809248
$ws.Range("F6").Value = 0.04046242774566474
$ws.Range("J6").Value = 0.2312138728323699
$ws.Range("O6").Value = 0.01734104046242774
$ws.Range("Q6").Value = 0.2023121387283237
$ws.Range("R6").Value = 0.06936416184971098
$ws.Range("S6").Value = 0.3872832369942196
$ws.Range("B7").Value = 0.1206030150753769
$ws.Range("D7").Value = 0.01507537688442211
$ws.Range("E7").Value = 0.005025125628140704
$ws.Range("F7").Value = 0.04020100502512563
$ws.Range("J7").Value = 0.1557788944723618
$ws.Range("O7").Value = 0.02512562814070352
$ws.Range("Q7").Value = 0.185929648241206
$ws.Range("R7").Value = 0.07537688442211055
$ws.Range("S7").Value = 0.3768844221105528
$ws.Range("B8").Value = 0.1098901098901099
$ws.Range("D8").Value = 0.02857142857142857
$ws.Range("F8").Value = 0.05494505494505494
$ws.Range("J8").Value = 0.1076923076923077
$ws.Range("O8").Value = 0.01318681318681319
$ws.Range("Q8").Value = 0.189010989010989
$ws.Range("R8").Value = 0.08351648351648351
$ws.Range("S8").Value = 0.4131868131868132
$ws.Range("B9").Value = 0.1104651162790698
$ws.Range("D9").Value = 0.01162790697674419
$ws.Range("F9").Value = 0.05232558139534884
$ws.Range("J9").Value = 0.08139534883720931
$ws.Range("O9").Value = 0.02325581395348837
$ws.Range("Q9").Value = 0.1453488372093023
$ws.Range("R9").Value = 0.1337209302325581
$ws.Range("S9").Value = 0.4418604651162791
$ws.Range("B10").Value = 0.1312292358803987
$ws.Range("D10").Value = 0.02408637873754153
$ws.Range("F10").Value = 0.06146179401993355
$ws.Range("J10").Value = 0.1104651162790698
$ws.Range("O10").Value = 0.02159468438538206
$ws.Range("Q10").Value = 0.1669435215946844
$ws.Range("R10").Value = 0.08803986710963455
$ws.Range("S10").Value = 0.3961794019933555
$ws.Range("G11").Value = 0.1538461538461539
$ws.Range("J11").Value = 0.1225071225071225
$ws.Range("K11").Value = 0.2393162393162393
$ws.Range("L11").Value = 0.4729344729344729
$ws.Range("S11").Value = 0.0113960113960114
$ws.Range("G12").Value = 0.6878612716763006
$ws.Range("J12").Value = 0.2543352601156069
$ws.Range("K12").Value = 0.005780346820809248
$ws.Range("L12").Value = 0.02890173410404624
$ws.Range("S12").Value = 0.02312138728323699
$ws.Range("G13").Value = 0.7142857142857143
$ws.Range("J13").Value = 0.2653061224489796
$ws.Range("S13").Value = 0.02040816326530612
$ws.Range("F15").Value = 0.02415458937198068
$ws.Range("H15").Value = 0.1400966183574879
$ws.Range("I15").Value = 0.05797101449275362
$ws.Range("J15").Value = 0.357487922705314
$ws.Range("K15").Value = 0.05314009661835749
$ws.Range("M15").Value = 0.01449275362318841
$ws.Range("O15").Value = 0.02898550724637681
$ws.Range("S15").Value = 0.3236714975845411
$ws.Range("H16").Value = 0.1818181818181818
$ws.Range("I16").Value = 0.08133971291866028
$ws.Range("J16").Value = 0.3684210526315789
$ws.Range("K16").Value = 0.1196172248803828
$ws.Range("M16").Value = 0.03827751196172249
$ws.Range("O16").Value = 0.02870813397129187
$ws.Range("S16").Value = 0.1818181818181818
$ws.Range("F17").Value = 0.0131578947368421
$ws.Range("H17").Value = 0.1710526315789474
$ws.Range("I17").Value = 0.07894736842105263
$ws.Range("J17").Value = 0.4263157894736842
$ws.Range("K17").Value = 0.09736842105263158
$ws.Range("M17").Value = 0.0131578947368421
$ws.Range("O17").Value = 0.06315789473684211
$ws.Range("S17").Value = 0.1368421052631579
$ws.Range("F18").Value = 0.0155440414507772
$ws.Range("H18").Value = 0.2435233160621762
$ws.Range("I18").Value = 0.1088082901554404
$ws.Range("J18").Value = 0.3523316062176166
$ws.Range("K18").Value = 0.06735751295336788
$ws.Range("M18").Value = 0.01036269430051814
$ws.Range("O18").Value = 0.07772020725388601
$ws.Range("S18").Value = 0.1243523316062176
$ws.Range("F19").Value = 0.01067073170731707
$ws.Range("H19").Value = 0.2126524390243902
$ws.Range("I19").Value = 0.07088414634146341
$ws.Range("J19").Value = 0.3521341463414634
$ws.Range("K19").Value = 0.1356707317073171
$ws.Range("M19").Value = 0.02591463414634146
$ws.Range("O19").Value = 0.06326219512195122
$ws.Range("S19").Value = 0.1288109756097561
